$d = $word.ActiveDocument

# Locate the word "dotnetcore".
$r = $d.Content
$r.Find.Execute("dotnetcore") | Out-Null
$start = $r.Start
$end = $r.End

# The word is currently one run: "dotnetcore" (rsidR="00C06B94").
# Target layout is three runs sharing identical character formatting:
#   "dot"    -> the original run, left untouched (keeps its rsid)
#   "N"      -> brand new run (the capitalised letter)
#   "etcore" -> brand new run (the remainder, re-typed so it loses the old rsid)

$nStart = $start + 3
$nEnd = $nStart + 1
$tailStart = $nEnd
$tailLen = $end - $tailStart

# Step 1: turn the lower-case "n" into "N".
# Toggling Bold around a genuine interior text change forces the engine to
# split the run instead of silently re-merging it, while leaving the
# untouched "dot" prefix and "etcore" suffix attached to the original run
# (so they keep its rsid) for now.
$rN = $d.Range($nStart, $nEnd)
$rN.Bold = 1
$rN.Text = "N"
$rN.Bold = 0

# After step 1, "etcore" is its own run (tailStart .. end) but it still
# carries the original rsid, since its text didn't actually change. Step 2:
# force a genuine (temporary) content change on it so the engine rebuilds it
# as a fresh run with no rsid, matching the target XML, then trim the extra
# character back off.
$tailText = $d.Range($tailStart, $end).Text
$rTail = $d.Range($tailStart, $end)
$rTail.Bold = 1
$rTail.Text = $tailText + "x"
$rExtra = $d.Range($tailStart + $tailLen, $tailStart + $tailLen + 1)
$rExtra.Text = ""
$rTailClean = $d.Range($tailStart, $tailStart + $tailLen)
$rTailClean.Bold = 0
